# Insert a new data row at row 252 (pushing existing rows 252-366 down to 253-367)
# and populate it with a new weekly price sample, as reflected in the commit
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 252..366 down to 253..367, making room for a new row 252.
$ws.Rows.Item(252).Insert()

# Populate the newly inserted row 252 with the new sample data.
$ws.Cells.Item(252, 1).Value = 3
$ws.Cells.Item(252, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(252, 3).Value = "Coquimbo"
$ws.Cells.Item(252, 4).Value = 44726
$ws.Cells.Item(252, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(252, 5).Value = 5
$ws.Cells.Item(252, 6).Value = 100112043
$ws.Cells.Item(252, 7).Value = "Pepino ensalada"
$ws.Cells.Item(252, 8).Value = "Sin especificar"
$ws.Cells.Item(252, 9).Value = "Primera"
$ws.Cells.Item(252, 10).Value = 145
$ws.Cells.Item(252, 11).Value = 18000
$ws.Cells.Item(252, 12).Value = 19000
$ws.Cells.Item(252, 13).Value = 18483
$ws.Cells.Item(252, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(252, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(252, 16).Value = 264
$ws.Cells.Item(252, 17).Value = 70
$ws.Cells.Item(252, 18).Value = "Hortaliza"
